$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.416.16"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.05%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.848.98"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.21%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.69%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6284"
$ws.Range("D6").Style = "Normal"

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07693"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.97%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2922"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.23%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.05"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.92%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07748"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.47%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.857.99"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.07%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.038"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.74%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.00001084"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.34%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6826"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.37%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.62"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.34%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.189"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.32%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.440.76"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.04%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "228.90"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.05%  "

$ws.Range("E20").Value = "  -0.16%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.04%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.463"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.12%  "

$ws.Range("E23").Value = "  +0.01%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "157.69"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.50%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1380"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.01%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.418"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.71%  "

$ws.Range("E27").Value = "  +0.78%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.354"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.05%  "

$ws.Range("E29").Value = "  +0.26%  "

$ws.Range("E30").Value = "  +0.08%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.125"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.58%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.048"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.64%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.845"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.07%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.165"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.62%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7083"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.30%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.594"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.09%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.225.87"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.72%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01791"
$ws.Range("D38").Style = "Normal"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.753"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.56%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.450"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.14%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9053"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.42%  "

$ws.Range("E42").Value = "  +0.02%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.027.42"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.82%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.87"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.13%  "

$ws.Range("E45").Value = "  +0.38%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.199"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.44%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000120"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.50%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4023"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.57%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1157"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.08%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.003"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.88%  "

$ws.Range("E51").Value = "  +0.25%  "
